$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filters save games) - updated B:E and G (sum) columns
$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 22898927661.19195
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 22898927866.10947

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 17.08608867836142

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 5.582307763322248

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 18.71679738969934
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 24.14949828602258

$ws.Range("B6").Value = 0.6545652718822623
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 0.7210945179870265
$ws.Range("E6").Value = 13.86384647080068
$ws.Range("G6").Value = 16.86649396021207

$ws.Range("B7").Value = 0.1169995834814548
$ws.Range("C7").Value = 0.3048912486333797
$ws.Range("D7").Value = 0.1496068669990043
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 1.104883657715537
